# Fixed naive component forecaster bug - Presentation state 11.02.
#
# For every data row (rows 2-24), a new "Q(-1)" style error value is
# inserted into column B. All of that row's existing values shift one
# column to the right, and anything that would overflow past column K
# (the 10th data column, B..K) is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to insert at column B for each row (row number => value)
$newValues = @{
    2  = 2.057869132359739
    3  = 6.652313087672924
    4  = -18.36749132628568
    5  = 7.513167073507937
    6  = 0.9564081874156993
    7  = -4.157449276732949
    8  = 1.546611864454844
    9  = 1.156631887942306
    10 = -1.025188112727922
    11 = 0.08364543516793629
    12 = -0.1538585523806955
    13 = 0.7495351060200912
    14 = 0.03849281619118239
    15 = -0.2590580299438133
    16 = 0.01855976243503714
    17 = 0.1467044301255134
    18 = -0.1819613811903656
    19 = 0.4718454808444464
    20 = -0.08594117411414147
    21 = -0.07695400962807622
    22 = -0.5068991247689255
    23 = 0.6215838649243215
    24 = -0.2766911554241067
}

# Column letters available for the data (B through K = 10 columns)
$cols = @("A","B","C","D","E","F","G","H","I","J","K")
$maxDataCol = 11   # column K

for ($row = 2; $row -le 24; $row++) {

    # Find the last populated column in columns B..K for this row.
    $lastColIndex = 1   # 1 = column A; start search from column B (index 2)
    for ($c = 2; $c -le $maxDataCol; $c++) {
        $cell = $ws.Cells.Item($row, $c)
        if ($cell.Value2 -ne $null) {
            $lastColIndex = $c
        }
    }

    $newVal = $newValues[$row]

    if ($lastColIndex -ge 2) {
        # Row currently has values in B..$lastColIndex; read them.
        $oldCount = $lastColIndex - 1   # number of existing values (B.. lastCol)
        $srcRange = $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, $lastColIndex))
        $oldVals = $srcRange.Value2

        # Build the shifted list: new value first, then old values, truncated
        # so it never extends past column K.
        $shifted = New-Object System.Collections.ArrayList
        [void]$shifted.Add($newVal)
        if ($oldCount -eq 1) {
            [void]$shifted.Add($oldVals)
        } else {
            for ($i = 1; $i -le $oldCount; $i++) {
                [void]$shifted.Add($oldVals[1, $i])
            }
        }

        $writeCount = [Math]::Min($shifted.Count, $maxDataCol - 1)
        $destLastCol = 1 + $writeCount

        $newArr = New-Object 'object[,]' 1, $writeCount
        for ($i = 0; $i -lt $writeCount; $i++) {
            $newArr[0, $i] = $shifted[$i]
        }

        $destRange = $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, $destLastCol))
        $destRange.Value2 = $newArr

        # If the row previously reached column K and nothing was dropped
        # (i.e. it didn't overflow), there's nothing further to clear.
        # If the row grew by one column compared to before, no stale cell
        # is left behind. Nothing else to clean up here because the shift
        # only ever keeps the same or fewer cells than before+1.
    } else {
        # Row had no values at all (e.g. row 24) - just write the new value.
        $ws.Cells.Item($row, 2).Value2 = $newVal
    }
}
